$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.097357333333333
$ws.Range("H2").Value = 12.292072
$ws.Range("I2").Value = 0.03407839216891784
$ws.Range("J2").Value = 0.03407839216891784
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 24.519512
$ws.Range("N2").Value = 73.558536
$ws.Range("O2").Value = 0.4736537296697991
$ws.Range("P2").Value = 0.4736537296697991
$ws.Range("Q2").Value = 100.4652023029547
$ws.Range("R2").Value = 904.186820726592
$ws.Range("S2").Value = 0.01614135755195801
$ws.Range("T2").Value = 0.01614135755195801

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.097357333333333
$ws.Range("H3").Value = 12.292072
$ws.Range("I3").Value = 0.03407839216891784
$ws.Range("J3").Value = 0.03407839216891784
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 14.70328633333333
$ws.Range("N3").Value = 44.109859
$ws.Range("O3").Value = 0.284029568377475
$ws.Range("P3").Value = 0.284029568377475
$ws.Range("Q3").Value = 60.24461808198311
$ws.Range("R3").Value = 542.201562737848
$ws.Range("S3").Value = 0.009679271018736058
$ws.Range("T3").Value = 0.009679271018736058

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.097357333333333
$ws.Range("H4").Value = 12.292072
$ws.Range("I4").Value = 0.03407839216891784
$ws.Range("J4").Value = 0.03407839216891784
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.543947
$ws.Range("N4").Value = 37.631841
$ws.Range("O4").Value = 0.2423167019527259
$ws.Range("P4").Value = 0.2423167019527259
$ws.Range("Q4").Value = 51.39703322939467
$ws.Range("R4").Value = 462.573299064552
$ws.Range("S4").Value = 0.008257763598223774
$ws.Range("T4").Value = 0.008257763598223774

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 114.2138366666667
$ws.Range("H5").Value = 342.64151
$ws.Range("I5").Value = 0.9499351900257489
$ws.Range("J5").Value = 0.9499351900257488
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 24.519512
$ws.Range("N5").Value = 73.558536
$ws.Range("O5").Value = 0.4736537296697991
$ws.Range("P5").Value = 0.4736537296697991
$ws.Range("Q5").Value = 2800.467538714373
$ws.Range("R5").Value = 25204.20784842936
$ws.Range("S5").Value = 0.4499403457002853
$ws.Range("T5").Value = 0.4499403457002852

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 114.2138366666667
$ws.Range("H6").Value = 342.64151
$ws.Range("I6").Value = 0.9499351900257489
$ws.Range("J6").Value = 0.9499351900257488
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 14.70328633333333
$ws.Range("N6").Value = 44.109859
$ws.Range("O6").Value = 0.284029568377475
$ws.Range("P6").Value = 0.284029568377475
$ws.Range("Q6").Value = 1679.318743738565
$ws.Range("R6").Value = 15113.86869364709
$ws.Range("S6").Value = 0.2698096820095881
$ws.Range("T6").Value = 0.269809682009588

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 114.2138366666667
$ws.Range("H7").Value = 342.64151
$ws.Range("I7").Value = 0.9499351900257489
$ws.Range("J7").Value = 0.9499351900257488
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 12.543947
$ws.Range("N7").Value = 37.631841
$ws.Range("O7").Value = 0.2423167019527259
$ws.Range("P7").Value = 0.2423167019527259
$ws.Range("Q7").Value = 1432.692313813323
$ws.Range("R7").Value = 12894.23082431991
$ws.Range("S7").Value = 0.2301851623158755
$ws.Range("T7").Value = 0.2301851623158754

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.922099666666667
$ws.Range("H8").Value = 5.766299
$ws.Range("I8").Value = 0.01598641780533329
$ws.Range("J8").Value = 0.01598641780533329
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 24.519512
$ws.Range("N8").Value = 73.558536
$ws.Range("O8").Value = 0.4736537296697991
$ws.Range("P8").Value = 0.4736537296697991
$ws.Range("Q8").Value = 47.12894584202934
$ws.Range("R8").Value = 424.160512578264
$ws.Range("S8").Value = 0.007572026417555798
$ws.Range("T8").Value = 0.007572026417555798

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.922099666666667
$ws.Range("H9").Value = 5.766299
$ws.Range("I9").Value = 0.01598641780533329
$ws.Range("J9").Value = 0.01598641780533329
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 14.70328633333333
$ws.Range("N9").Value = 44.109859
$ws.Range("O9").Value = 0.284029568377475
$ws.Range("P9").Value = 0.284029568377475
$ws.Range("Q9").Value = 28.26118176020456
$ws.Range("R9").Value = 254.350635841841
$ws.Range("S9").Value = 0.004540615349150796
$ws.Range("T9").Value = 0.004540615349150796

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.922099666666667
$ws.Range("H10").Value = 5.766299
$ws.Range("I10").Value = 0.01598641780533329
$ws.Range("J10").Value = 0.01598641780533329
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 12.543947
$ws.Range("N10").Value = 37.631841
$ws.Range("O10").Value = 0.2423167019527259
$ws.Range("P10").Value = 0.2423167019527259
$ws.Range("Q10").Value = 24.11071634738434
$ws.Range("R10").Value = 216.996447126459
$ws.Range("S10").Value = 0.003873776038626698
$ws.Range("T10").Value = 0.003873776038626698
